$wb = $excel.ActiveWorkbook

# Sheet "展览" (1st sheet) - column F updates
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F4").Value = 272
$ws1.Range("F6").Value = 10141
$ws1.Range("F9").Value = 1261
$ws1.Range("F10").Value = 6592
$ws1.Range("F12").Value = 419
$ws1.Range("F15").Value = 3119
$ws1.Range("F20").Value = 26
$ws1.Range("F23").Value = 1557

# Sheet "全部类型" (4th sheet) - column F updates
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F5").Value = 272
$ws4.Range("F7").Value = 10141
$ws4.Range("F10").Value = 1261
$ws4.Range("F11").Value = 6592
$ws4.Range("F13").Value = 419
$ws4.Range("F16").Value = 3119
$ws4.Range("F21").Value = 26
$ws4.Range("F24").Value = 1557
